# Auto-generated Excel COM-interop script applying the Typhon_Profits.xlsx diff.
# For each affected row/column: sets new numeric values, clears cells that were
# removed in the diff, and adds values to cells that were newly introduced.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H58").Value = 1665.7693
$ws.Range("I58").Value = 231
$ws.Range("J58").Value = 2562.5
$ws.Range("K58").Value = 693
$ws.Range("L58").Value = 7687.5
$ws.Range("M58").Value = -543
$ws.Range("N58").Value = -7987.5
$ws.Range("H62").Value = 2088.2856
$ws.Range("I62").Value = 1632.5
$ws.Range("J62").Value = 2999.8572
$ws.Range("K62").Value = 1632.5
$ws.Range("L62").Value = 2999.8572
$ws.Range("M62").Value = -1008.5
$ws.Range("N62").Value = -4247.8572
$ws.Range("H65").Value = 2088.2856
$ws.Range("I65").Value = 1632.5
$ws.Range("J65").Value = 2999.8572
$ws.Range("K65").Value = 8162.5
$ws.Range("L65").Value = 14999.286
$ws.Range("M65").Value = -5042.5
$ws.Range("N65").Value = -21239.286
$ws.Range("H107").Value = 876.38464
$ws.Range("I107").Value = 784.8570999999999
$ws.Range("J107").Value = 983.1667
$ws.Range("K107").Value = 784.8570999999999
$ws.Range("L107").Value = 983.1667
$ws.Range("M107").Value = 1135.1429
$ws.Range("N107").Value = -4823.1667
$ws.Range("H129").Value = 849.6087
$ws.Range("J129").Value = 849.5909
$ws.Range("L129").Value = 2548.7727
$ws.Range("N129").Value = -12548.7727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 25532.137
$ws.Range("I132").Value = 2619.182
$ws.Range("K132").Value = 7857.545999999999
$ws.Range("M132").Value = -5327.545999999999
$ws.Range("H139").Value = 48000
$ws.Range("J139").Value = 48000
$ws.Range("L139").Value = 48000
$ws.Range("N139").Value = -58280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 19830.285
$ws.Range("J81").Value = 19830.285
$ws.Range("L81").Value = 19830.285
$ws.Range("N81").Value = -21952.285
$ws.Range("H84").Value = 19830.285
$ws.Range("J84").Value = 19830.285
$ws.Range("L84").Value = 59490.855
$ws.Range("N84").Value = -70098.855
$ws.Range("H86").Value = 1935.3846
$ws.Range("I86").Value = 1781
$ws.Range("J86").Value = 2450
$ws.Range("K86").Value = 1781
$ws.Range("L86").Value = 2450
$ws.Range("M86").Value = -658
$ws.Range("N86").Value = -4696
$ws.Range("H89").Value = 1935.3846
$ws.Range("I89").Value = 1781
$ws.Range("J89").Value = 2450
$ws.Range("K89").Value = 8905
$ws.Range("L89").Value = 12250
$ws.Range("M89").Value = -3289
$ws.Range("N89").Value = -23482
$ws.Range("H134").Value = 8845.25
$ws.Range("I134").Value = 9680.286
$ws.Range("K134").Value = 29040.858
$ws.Range("M134").Value = -26505.858
$ws.Range("H135").Value = 46744
$ws.Range("J135").Value = 46744
$ws.Range("L135").Value = 46744
$ws.Range("N135").Value = -56884

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 3999.8
$ws.Range("J14").Value = 3000
$ws.Range("L14").Value = 3000
$ws.Range("N14").Value = -3340
$ws.Range("H62").Value = 6178.3335
$ws.Range("I62").Value = 4026.25
$ws.Range("K62").Value = 4026.25
$ws.Range("M62").Value = -3402.25
$ws.Range("H65").Value = 6178.3335
$ws.Range("I65").Value = 4026.25
$ws.Range("K65").Value = 20131.25
$ws.Range("M65").Value = -17011.25
$ws.Range("H99").Value = 4772.6665
$ws.Range("I99").Value = 3514.1333
$ws.Range("K99").Value = 3514.1333
$ws.Range("M99").Value = -2016.1333
$ws.Range("H105").Value = 11364342
$ws.Range("I105").Value = 17857620
$ws.Range("J105").Value = 1105.5
$ws.Range("K105").Value = 17857620
$ws.Range("L105").Value = 1105.5
$ws.Range("M105").Value = -17855873
$ws.Range("N105").Value = -4599.5
$ws.Range("H126").Value = 4772.6665
$ws.Range("I126").Value = 3514.1333
$ws.Range("K126").Value = 10542.3999
$ws.Range("M126").Value = -8072.3999
$ws.Range("H132").Value = 24233.043
$ws.Range("I132").Value = 28907.5
$ws.Range("J132").Value = 7405
$ws.Range("K132").Value = 86722.5
$ws.Range("L132").Value = 22215
$ws.Range("M132").Value = -84192.5
$ws.Range("N132").Value = -27275

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("H44").Value = 170.83333
$ws.Range("J44").Value = 241.66667
$ws.Range("L44").Value = 725.00001
$ws.Range("N44").Value = -1521.00001
$ws.Range("H51").Value = 2519.0908
$ws.Range("J51").Value = 2634.4443
$ws.Range("L51").Value = 7903.3329
$ws.Range("N51").Value = -8823.332900000001
$ws.Range("H55").Value = 2478.2
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 2478.2
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 7434.599999999999
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -7788.599999999999
$ws.Range("H81").Value = 4237.75
$ws.Range("I81").Value = 3000
$ws.Range("J81").Value = 4320.2666
$ws.Range("K81").Value = 9000
$ws.Range("L81").Value = 12960.7998
$ws.Range("M81").Value = -7877
$ws.Range("N81").Value = -15206.7998
$ws.Range("H84").Value = 4237.75
$ws.Range("I84").Value = 3000
$ws.Range("J84").Value = 4320.2666
$ws.Range("K84").Value = 27000
$ws.Range("L84").Value = 38882.3994
$ws.Range("M84").Value = -21384
$ws.Range("N84").Value = -50114.3994
$ws.Range("H88").Value = 9863.200000000001
$ws.Range("J88").Value = 9863.200000000001
$ws.Range("L88").Value = 29589.6
$ws.Range("N88").Value = -30445.6
$ws.Range("H91").Value = 9863.200000000001
$ws.Range("J91").Value = 9863.200000000001
$ws.Range("L91").Value = 29589.6
$ws.Range("N91").Value = -32553.6
$ws.Range("H94").Value = 4259.9165
$ws.Range("I94").Value = 825.5
$ws.Range("J94").Value = 4946.8
$ws.Range("K94").Value = 2476.5
$ws.Range("L94").Value = 14840.4
$ws.Range("M94").Value = -1800.5
$ws.Range("N94").Value = -16192.4
$ws.Range("H113").Value = 447.29413
$ws.Range("I113").Value = 429.33334
$ws.Range("J113").Value = 467.5
$ws.Range("K113").Value = 1288.00002
$ws.Range("L113").Value = 1402.5
$ws.Range("M113").Value = 881.9999800000001
$ws.Range("N113").Value = -5742.5
$ws.Range("H114").Value = 165
$ws.Range("I114").Value = 165
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 495
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = 2759
$ws.Range("N114").ClearContents()
$ws.Range("H116").Value = 500
$ws.Range("I116").Value = 500
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1500
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1942
$ws.Range("N116").ClearContents()
$ws.Range("H123").Value = 3404.2
$ws.Range("I123").Value = 1010.3333
$ws.Range("J123").Value = 6995
$ws.Range("K123").Value = 3030.9999
$ws.Range("L123").Value = 20985
$ws.Range("M123").Value = -580.9998999999998
$ws.Range("N123").Value = -25885
$ws.Range("H131").Value = 724.71
$ws.Range("I131").Value = 302.5
$ws.Range("J131").Value = 742.30206
$ws.Range("K131").Value = 907.5
$ws.Range("L131").Value = 2226.90618
$ws.Range("M131").Value = 4132.5
$ws.Range("N131").Value = -12306.90618

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1674.6
$ws.Range("I122").Value = 1370
$ws.Range("K122").Value = 4110
$ws.Range("M122").Value = -1660
$ws.Range("H132").Value = 25540.436
$ws.Range("I132").Value = 4079.4614
$ws.Range("K132").Value = 12238.3842
$ws.Range("M132").Value = -9708.3842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3297.6
$ws.Range("I40").Value = 2762
$ws.Range("K40").Value = 2762
$ws.Range("M40").Value = -2626
$ws.Range("H100").Value = 2714.9048
$ws.Range("I100").Value = 1863
$ws.Range("J100").Value = 3652
$ws.Range("K100").Value = 1863
$ws.Range("L100").Value = 3652
$ws.Range("M100").Value = -1322
$ws.Range("N100").Value = -4734
$ws.Range("H132").Value = 3922.4
$ws.Range("I132").Value = 1304
$ws.Range("J132").Value = 4577
$ws.Range("K132").Value = 3912
$ws.Range("L132").Value = 13731
$ws.Range("M132").Value = -1382
$ws.Range("N132").Value = -18791

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2674780.5
$ws.Range("I107").Value = 282.4
$ws.Range("J107").Value = 6495492
$ws.Range("K107").Value = 847.1999999999999
$ws.Range("L107").Value = 19486476
$ws.Range("M107").Value = 1072.8
$ws.Range("N107").Value = -19490316
$ws.Range("H132").Value = 4070.8572
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 4082.6667
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 12248.0001
$ws.Range("M132").Value = -9470
$ws.Range("N132").Value = -17308.0001
$ws.Range("H136").Value = 1275.5294
$ws.Range("I136").Value = 971.13635
$ws.Range("K136").Value = 2913.40905
$ws.Range("M136").Value = -363.4090500000002
